$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '89.423.02'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.066.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.94%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +8.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '616.18'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -7.13%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.070.26'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.706'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.99%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.73'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.207.30'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.34'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.637.41'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.077.54'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.70'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000214'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.68'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -6.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '429.18'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -8.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.37'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.64'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.52'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.97'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -10.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.59'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.238.96'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.46%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.98'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.155'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.36%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -12.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.49'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.149'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.99'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +64.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.97'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '487.95'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.59'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.24'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0895'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.07'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.395'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -8.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.22'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.21%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.669'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -11.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.19'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.29'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.10%  '
